# Update "想去人数" (interest count) figures in the F column, sourced from a
# refreshed data pull (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 325
$ws1.Range("F9").Value  = 300
$ws1.Range("F10").Value = 1684
$ws1.Range("F11").Value = 336
$ws1.Range("F12").Value = 1392
$ws1.Range("F13").Value = 788
$ws1.Range("F14").Value = 319
$ws1.Range("F15").Value = 659
$ws1.Range("F16").Value = 12635
$ws1.Range("F17").Value = 12659
$ws1.Range("F18").Value = 934
$ws1.Range("F23").Value = 500
$ws1.Range("F24").Value = 1979
$ws1.Range("F25").Value = 22
$ws1.Range("F27").Value = 230

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 70
$ws2.Range("F11").Value = 5

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 151

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 151
$ws4.Range("F12").Value = 325
$ws4.Range("F14").Value = 300
$ws4.Range("F15").Value = 1684
$ws4.Range("F16").Value = 336
$ws4.Range("F17").Value = 1392
$ws4.Range("F18").Value = 788
$ws4.Range("F19").Value = 319
$ws4.Range("F20").Value = 70
$ws4.Range("F21").Value = 659
$ws4.Range("F22").Value = 12635
$ws4.Range("F23").Value = 12659
$ws4.Range("F24").Value = 934
$ws4.Range("F29").Value = 500
$ws4.Range("F32").Value = 1979
$ws4.Range("F33").Value = 22
$ws4.Range("F37").Value = 230
$ws4.Range("F40").Value = 5
